# Update the "想去人数" (F column) figures for several events on both the
# "展览" and "全部类型" worksheets, reflecting refreshed stats from the
# site regeneration (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 246
$ws1.Range("F3").Value = 77
$ws1.Range("F4").Value = 795
$ws1.Range("F5").Value = 516

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 246
$ws4.Range("F3").Value = 77
$ws4.Range("F4").Value = 795
$ws4.Range("F6").Value = 516
